$d = $word.ActiveDocument

# --- Edit 1: merge the "September 10 - 16" paragraph's three runs (which
# were split around a proofed "as") back into a single run / single w:t,
# dropping the proofErr gramStart/gramEnd markers in the process. A
# Find/Replace across the whole sentence collapses the run back to one. ---
$oldText1 = "September 10 " + [char]0x2013 + " 16: Worked on the prototype of the game. Drew up wireframes of how the game is supposed to look as a bare bones structure. Edited the design document when I was finally able to join a group and see the submission requirements."
$find1 = $d.Content.Find
$find1.ClearFormatting()
$find1.Text = $oldText1
$find1.Replacement.ClearFormatting()
$find1.Replacement.Text = $oldText1
$find1.Execute($find1.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2)

# --- Edit 2: append a new blank paragraph followed by a new "Nov 5 - Nov 20"
# timelog paragraph at the end of the document (before the sectPr). ---
$newText2 = "Nov 5 " + [char]0x2013 + " Nov 20: Fleshed out the simulation questions and have a general idea of how to upscale in difficulty if time permits. Also got the hang of persistent data and I am working on typing up the code into the sim file so organization of the scenes is maintained."

$lastPara = $d.Paragraphs.Last
$endRange = $lastPara.Range
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

$endRange2 = $d.Paragraphs.Last.Range
$endRange2.Collapse(0)
$endRange2.InsertParagraphAfter()

$finalRange = $d.Content
$finalRange.Collapse(0)
$finalRange.InsertAfter($newText2)

Write-Output "ParagraphCount: $($d.Paragraphs.Count)"
